# 2024_전기공사_통합데이터.xlsx edit:
# Column K ("투찰률") is removed. For rows 16-115, the numeric value that
# used to live in column K is moved into column H ("낙찰률"), which was
# previously blank for those rows. Rows 1-15 already had H populated and
# had only an empty placeholder in K, so nothing needs to move there -
# deleting the column takes care of them. Finally the sheet's used
# dimension shrinks from A1:K115 to A1:J115, which Excel recomputes
# automatically once the column is deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column 11 = K, column 8 = H.
for ($r = 16; $r -le 115; $r++) {
    $kCell = $ws.Cells.Item($r, 11)
    $kVal = $kCell.Value2
    if ($kVal -ne $null -and $kVal -ne "") {
        $ws.Cells.Item($r, 8).Value = $kVal
    }
}

# Remove column K entirely (header + data), shifting the dimension to A1:J115.
$ws.Columns.Item(11).Delete()
